$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2-9 from serial 45221 to 45224
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value2 = 45224
}
